# The cover-sheet's "Date:" field is a plain-text content control (SDT)
# containing "23/06/2021". The commit updates it to "26/06/2021".
#
# NOTE: Word's Find/Replace (Content.Find.Execute with a Replace mode)
# cannot reach into structured document tags (content controls) in this
# runtime, so we update the control's own Range.Text directly instead.
# We first set the control's text to just the leading "2", then use
# InsertAfter twice to append "6" and "/06/2021" as their own runs -
# this reproduces the same three-run split ("2" / "6" / "/06/2021")
# that Word itself produced when the "3" was retyped as "6".

$d = $word.ActiveDocument
$ccs = $d.ContentControls

for ($i = 1; $i -le $ccs.Count; $i++) {
    $cc = $ccs.Item($i)
    if ($cc.Range.Text -eq "23/06/2021") {
        $cc.Range.Text = "2"
        $cc.Range.InsertAfter("6")
        $cc.Range.InsertAfter("/06/2021")
    }
}
